# Update the cryptos price/volume table with refreshed values.
# Price values that look numeric (e.g. "1.014") are written with a
# leading apostrophe so Excel stores them as text (matching the
# original inline-string cell type) instead of coercing them to
# numbers and losing formatting (trailing zeros, scientific
# notation, etc). Values that already contain multiple dots
# (e.g. "28.948.07") are never auto-coerced by Excel, so they are
# assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.948.07'
$ws.Range("E2").Value = '  -2.26%  '

$ws.Range("D3").Value = '1.975.41'
$ws.Range("E3").Value = '  -1.63%  '

$ws.Range("D4").Value = '''1.014'
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = '''326.90'
$ws.Range("E5").Value = '  -1.50%  '

$ws.Range("D6").Value = '''1.012'
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D7").Value = '''0.4845'
$ws.Range("E7").Value = '  -3.88%  '

$ws.Range("D8").Value = '''0.4095'
$ws.Range("E8").Value = '  -3.90%  '

$ws.Range("D9").Value = '''54.35'
$ws.Range("E9").Value = '  -1.24%  '

$ws.Range("D10").Value = '''0.08708'
$ws.Range("E10").Value = '  -5.11%  '

$ws.Range("D11").Value = '''1.075'
$ws.Range("E11").Value = '  -4.59%  '

$ws.Range("D12").Value = '''22.56'
$ws.Range("E12").Value = '  -4.79%  '

$ws.Range("D13").Value = '1.939.53'
$ws.Range("E13").Value = '  -3.91%  '

$ws.Range("D14").Value = '''7.803'
$ws.Range("E14").Value = '  -4.03%  '

$ws.Range("D15").Value = '''6.310'
$ws.Range("E15").Value = '  -3.73%  '

$ws.Range("D16").Value = '''1.015'
$ws.Range("E16").Value = '  +0.17%  '

$ws.Range("D17").Value = '''91.26'
$ws.Range("E17").Value = '  -4.32%  '

$ws.Range("D18").Value = '''0.00001079'
$ws.Range("E18").Value = '  -4.30%  '

$ws.Range("D19").Value = '''0.06655'
$ws.Range("E19").Value = '  -0.39%  '

$ws.Range("D20").Value = '''18.97'
$ws.Range("E20").Value = '  -4.78%  '

$ws.Range("D21").Value = '''1.012'
$ws.Range("E21").Value = '  +0.16%  '

$ws.Range("D22").Value = '''5.828'
$ws.Range("E22").Value = '  -2.89%  '

$ws.Range("D23").Value = '28.955.76'
$ws.Range("E23").Value = '  -2.25%  '

$ws.Range("D24").Value = '''11.67'
$ws.Range("E24").Value = '  -3.27%  '

$ws.Range("D25").Value = '''2.293'
$ws.Range("E25").Value = '  +0.52%  '

$ws.Range("D26").Value = '2.170.46'
$ws.Range("E26").Value = '  -4.36%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''20.47'
$ws.Range("E27").Value = '  -1.78%  '

$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '''154.68'
$ws.Range("E28").Value = '  -2.98%  '

$ws.Range("D29").Value = '''6.102'
$ws.Range("E29").Value = '  -5.76%  '

$ws.Range("D30").Value = '''2.186'
$ws.Range("E30").Value = '  -6.56%  '

$ws.Range("D31").Value = '''125.54'
$ws.Range("E31").Value = '  -2.54%  '

$ws.Range("D32").Value = '''1.017'
$ws.Range("E32").Value = '  -4.60%  '

$ws.Range("D33").Value = '''0.09715'
$ws.Range("E33").Value = '  -2.58%  '

$ws.Range("D34").Value = '''1.485'
$ws.Range("E34").Value = '  -6.42%  '

$ws.Range("D35").Value = '''5.744'
$ws.Range("E35").Value = '  -2.22%  '

$ws.Range("D36").Value = '''3.711'
$ws.Range("E36").Value = '  -2.83%  '

$ws.Range("D37").Value = '''0.02375'
$ws.Range("E37").Value = '  -4.37%  '

$ws.Range("D38").Value = '''1.292'
$ws.Range("E38").Value = '  -2.43%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '''0.06303'
$ws.Range("E39").Value = '  -1.53%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''8.910'
$ws.Range("E40").Value = '  -7.04%  '

$ws.Range("D41").Value = '''0.6320'
$ws.Range("E41").Value = '  -4.27%  '

$ws.Range("D42").Value = '''11.27'
$ws.Range("E42").Value = '  -4.48%  '

$ws.Range("D43").Value = '''1.012'
$ws.Range("E43").Value = '  +0.14%  '

$ws.Range("D44").Value = '''0.1935'
$ws.Range("E44").Value = '  -6.98%  '

$ws.Range("D45").Value = '''1.355'
$ws.Range("E45").Value = '  +4.00%  '

$ws.Range("D46").Value = '''0.6042'
$ws.Range("E46").Value = '  -5.42%  '

$ws.Range("D47").Value = '''13.07'
$ws.Range("E47").Value = '  -4.28%  '

$ws.Range("D48").Value = '''2.113'
$ws.Range("E48").Value = '  -4.91%  '

$ws.Range("D49").Value = '''3.450'
$ws.Range("E49").Value = '  -2.27%  '

$ws.Range("D50").Value = '''0.00000000336'
$ws.Range("E50").Value = '  +3.74%  '

$ws.Range("D51").Value = '''2.175'
$ws.Range("E51").Value = '  +7.49%  '
